# Rename the "*_multiplier" input columns to "*_proportion" to better
# reflect what they represent, and add a bit more info / clarity for the
# readme-adjacent column headers. Underlying numeric/text data for every
# row is unchanged -- only the header labels below are updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q1").Value = "y_2025_proportion"
$ws.Range("R1").Value = "y_2026_proportion"
$ws.Range("S1").Value = "y_2027_proportion"
$ws.Range("T1").Value = "y_2028_proportion"
$ws.Range("U1").Value = "y_2029_proportion"
$ws.Range("V1").Value = "y_2030_proportion"
$ws.Range("X1").Value = "specialized_obesity_proportion"
$ws.Range("Z1").Value = "specialized_chemphys_proportion"
$ws.Range("AB1").Value = "specialized_lifestyle_proportion"

# Restore the cursor/selection to where the author last left it (K15)
# before saving.
$ws.Range("K15").Select()
